# Adds a new "2021" column (R) to the unemployment-rate table, mirroring
# the existing year columns (D:Q), and moves the active selection to S1.
#
# Column layout reminder: A=label(kk)/B=label(ru)/C=label(en), D..Q = years
# 2007..2020, and the new R = 2021. Row 4 holds the year headers, rows
# 5-43 hold the data (some rows are bold "total" rows, some are blank
# section headers, row 43 is the bottom, bordered row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new 2021 values, keyed by row -----------------------------------
$values = @{
  5  = 5.3
  6  = 6.3
  7  = 4.7
  9  = 6.6
  10 = 7.5
  11 = 6.2
  12 = 11.8
  13 = 15.5
  14 = 9.7
  15 = 6.3
  16 = 7.5
  17 = 5.6
  18 = 6.3
  19 = 10.8
  20 = 4.3
  21 = 1.9
  22 = 3.1
  23 = 1.1
  24 = 2.6
  25 = 3.8
  26 = 1.7
  27 = 5.3
  28 = 6.2
  29 = 4.8
  30 = 4.1
  31 = 3.3
  32 = 4.9
  33 = 2.8
  34 = 3.4
  35 = 2.6
  37 = 15.7
  38 = 7.9
  39 = 4.5
  40 = 4.4
  41 = 2.9
  42 = 1.4
}

# Rows whose label cell (column A) is bold ("Kyrgyz Republic" / regional
# / city totals) -- these get the bold, non-bordered, plain-number style
# instead of the regular one-decimal data style.
$boldRows = @(5, 9, 12, 15, 18, 21, 24, 27, 30, 33)

# Blank section-header rows ("by territory" / "by age group") - column R
# stays empty here too, just formatted like the rest of the row.
$blankRows = @(8, 36)

# Row 4 is the year header row; give it the same look as D4:Q4.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Cells.Item(4, 18).Value = 2021

foreach ($r in 5..42) {
  if ($boldRows -contains $r) {
    # Bold "total" rows: match the bold label-cell formatting (no border,
    # no decimal number format), same as column A/B/C on this row.
    $ws.Range("A$r").Copy()
    $ws.Range("R$r").PasteSpecial(-4122)
  } elseif ($blankRows -contains $r) {
    # Section-header rows: plain style, left empty.
    $ws.Range("G1").Copy()
    $ws.Range("R$r").PasteSpecial(-4122)
    continue
  } else {
    # Regular data rows: plain (non-bold) style, same as the rest of the
    # row's label cells.
    $ws.Range("C10").Copy()
    $ws.Range("R$r").PasteSpecial(-4122)
  }

  if ($values.ContainsKey($r)) {
    $ws.Cells.Item($r, 18).Value = $values[$r]
  }
}

# Row 43: bottom, bordered row. The 70-and-over age group has no 2021
# figure yet, so it gets the same "…" placeholder used elsewhere in that
# row (e.g. D43/E43/F43).
$ws.Range("D43").Copy()
$ws.Range("R43").PasteSpecial(-4122)
$ws.Cells.Item(43, 18).Value = "…"

# Move the selection the same way the author's workbook ended up (one
# column to the right of the new data, row 1).
$ws.Range("S1").Select()
